# Generate Report for Handoff
# Adds two new handed-off files (2fed40f3-... and a5d06b60-...) as new rows
# to the Overview / zh-cn / de-de sheets + tables, mirroring the existing
# "Ready for handoff" rows (e.g. 2c19810b-...).

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> table "Overview", rows 4 & 5
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Row 4
$wsOv.Range("A4").Value = "2fed40f3-050f-4370-ab23-effda179b47c.md"
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("E4").Value = "Ready for handoff"
$wsOv.Range("F4").Value = "Ready for handoff"
$wsOv.Range("G4").Value = "2016-08-13 10:50:01"
$wsOv.Range("G4").NumberFormat = $dateFmt
$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/94a0a7e9df8728f0ac05d4a0bd8c10e0ce681aff/e2e/2fed40f3-050f-4370-ab23-effda179b47c.md", "", "", "e2e\2fed40f3-050f-4370-ab23-effda179b47c.md")

# Row 5
$wsOv.Range("A5").Value = "a5d06b60-17cf-44f3-b757-6f4b919228e3.md"
$wsOv.Range("C5").Value = ".md"
$wsOv.Range("E5").Value = "Ready for handoff"
$wsOv.Range("F5").Value = "Ready for handoff"
$wsOv.Range("G5").Value = "2016-08-13 10:50:01"
$wsOv.Range("G5").NumberFormat = $dateFmt
$wsOv.Hyperlinks.Add($wsOv.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a932399d663578dfeb24042c2b9c88c4dfbc8329/e2e/a5d06b60-17cf-44f3-b757-6f4b919228e3.md", "", "", "e2e\a5d06b60-17cf-44f3-b757-6f4b919228e3.md")

# Expand the backing table + autofilter to cover the two new rows
$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G5"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> table "zh-cn", rows 4 & 5
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 4
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "2fed40f3-050f-4370-ab23-effda179b47c.bff8eea9fea42853be414b7a071ccd63607086eb.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-13 10:49:53"
$wsZh.Range("H4").NumberFormat = $dateFmt
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("K4").NumberFormat = $dateFmt
$wsZh.Range("M4").Value = "True"
$wsZh.Range("O4").Value = "False"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/94a0a7e9df8728f0ac05d4a0bd8c10e0ce681aff/e2e/2fed40f3-050f-4370-ab23-effda179b47c.md", "", "", "2fed40f3-050f-4370-ab23-effda179b47c.md")

# Row 5
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "a5d06b60-17cf-44f3-b757-6f4b919228e3.30e8ead1a3b0dc3f3d31e07b2151d916235d57a2.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-08-13 10:49:53"
$wsZh.Range("H5").NumberFormat = $dateFmt
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = $dateFmt
$wsZh.Range("M5").Value = "True"
$wsZh.Range("O5").Value = "False"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a932399d663578dfeb24042c2b9c88c4dfbc8329/e2e/a5d06b60-17cf-44f3-b757-6f4b919228e3.md", "", "", "a5d06b60-17cf-44f3-b757-6f4b919228e3.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> table "de-de", rows 4 & 5
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 4
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "2fed40f3-050f-4370-ab23-effda179b47c.bff8eea9fea42853be414b7a071ccd63607086eb.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-13 10:50:01"
$wsDe.Range("H4").NumberFormat = $dateFmt
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("K4").NumberFormat = $dateFmt
$wsDe.Range("M4").Value = "True"
$wsDe.Range("O4").Value = "False"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/94a0a7e9df8728f0ac05d4a0bd8c10e0ce681aff/e2e/2fed40f3-050f-4370-ab23-effda179b47c.md", "", "", "2fed40f3-050f-4370-ab23-effda179b47c.md")

# Row 5
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "a5d06b60-17cf-44f3-b757-6f4b919228e3.30e8ead1a3b0dc3f3d31e07b2151d916235d57a2.de-de.xlf"
$wsDe.Range("H5").Value = "2016-08-13 10:50:01"
$wsDe.Range("H5").NumberFormat = $dateFmt
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = $dateFmt
$wsDe.Range("M5").Value = "True"
$wsDe.Range("O5").Value = "False"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a932399d663578dfeb24042c2b9c88c4dfbc8329/e2e/a5d06b60-17cf-44f3-b757-6f4b919228e3.md", "", "", "a5d06b60-17cf-44f3-b757-6f4b919228e3.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))
